$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.669.81"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.69%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.588.22"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -2.53%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.53%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "207.04"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.06%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -3.32%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.56%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.21"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -4.64%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.71%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.79%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0866"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.64%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.814.24"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.53%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.588.40"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -2.36%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -4.09%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -4.77%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.652.77"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.78%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.43"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.37%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "219.10"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -4.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0696"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -3.03%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -3.89%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -4.86%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -3.46%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -3.74%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.62"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.12%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.10%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.69%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.36%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -4.64%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.21%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.71%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -5.42%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.371.39"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.21%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -5.71%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -4.71%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.979"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.65%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.90%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.535"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -3.38%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.824"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -3.41%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.54%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.973"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.82%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.49%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.57%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.20"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -4.09%  "
$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.725.37"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.53%  "
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.72"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -5.09%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.41"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.33%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.64%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0966"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -4.44%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0495"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.60%  "
